$d = $word.ActiveDocument

$replacements = @(
    @{old = "334×2="; new = "638×5="},
    @{old = "947×8="; new = "114×6="},
    @{old = "289×3="; new = "687×9="},
    @{old = "206×5="; new = "534×9="},
    @{old = "491×9="; new = "629×2="},
    @{old = "314×5="; new = "797×3="},
    @{old = "847×3="; new = "119×3="},
    @{old = "128×5="; new = "146×7="},
    @{old = "914×5="; new = "493×2="},
    @{old = "125×3="; new = "301×2="},
    @{old = "831×6="; new = "989×3="},
    @{old = "650×8="; new = "743×4="},
    @{old = "476×4="; new = "397×2="},
    @{old = "847×2="; new = "926×4="},
    @{old = "826×4="; new = "891×6="},
    @{old = "658×9="; new = "919×5="},
    @{old = "479×4="; new = "515×6="},
    @{old = "899×4="; new = "932×6="},
    @{old = "646×9="; new = "180×3="},
    @{old = "885×6="; new = "612×2="},
    @{old = "877×5="; new = "996×4="},
    @{old = "985×4="; new = "806×2="},
    @{old = "616×3="; new = "560×9="},
    @{old = "627×9="; new = "102×2="},
    @{old = "160×3="; new = "518×5="}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $r.new, 2)
}
